{"js": "// Update the date title and the 25 division problems (5x5 grid across\n// non-blank rows 0,4,8,12,16 of the single table), matching the document's\n// original top-to-bottom reading order. Several old/new strings repeat\n// across the list (e.g. \"48\u00f75=\" is both a replaced value and a replacement\n// target), so we must replace by POSITION, not by a global text search.\n\nconst titleMap = [\"2025-04-06 Sunday\"];\n\nconst cellRows = [\n  [\"14\u00f75=\", \"48\u00f75=\", \"59\u00f72=\", \"65\u00f76=\", \"79\u00f75=\"],\n  [\"47\u00f77=\", \"36\u00f79=\", \"16\u00f72=\", \"37\u00f77=\", \"74\u00f72=\"],\n  [\"40\u00f74=\", \"31\u00f76=\", \"41\u00f75=\", \"23\u00f74=\", \"93\u00f77=\"],\n  [\"96\u00f79=\", \"51\u00f73=\", \"46\u00f74=\", \"22\u00f74=\", \"23\u00f74=\"],\n  [\"56\u00f79=\", \"17\u00f78=\", \"16\u00f74=\", \"39\u00f75=\", \"90\u00f78=\"],\n];\n\nconst body = context.document.body;\n\n// --- Title paragraph (first paragraph in the body) ---\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nparas.items[0].getRange().insertText(titleMap[0], \"Replace\");\n\n// --- Table cells ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\nconst rowCount = rows.items.length;\n\n// Non-blank data rows are every 4th row starting at 0 (0,4,8,12,16).\nlet dataRowIdx = 0;\nfor (let r = 0; r < rowCount && dataRowIdx < cellRows.length; r += 4) {\n  const values = cellRows[dataRowIdx];\n  for (let c = 0; c < values.length; c++) {\n    const cell = table.getCell(r, c);\n    cell.getRange().insertText(values[c], \"Replace\");\n  }\n  dataRowIdx++;\n}\n\nawait context.sync();\n", "ps1": "# Update the date title and the 25 division problems (5x5 grid across\n# non-blank rows 1,5,9,13,17 of the single table), matching the document's\n# original top-to-bottom reading order. Several old/new strings repeat\n# across the list (e.g. \"48\u00f75=\" is both a replaced value and a replacement\n# target), so we replace by POSITION (row/column), not by a global\n# Find/Replace across the whole document.\n\n$d = $word.ActiveDocument\n\n# --- Title paragraph (first paragraph in the body) ---\n$d.Paragraphs.Item(1).Range.Text = \"2025-04-06 Sunday\"\n\n# --- Table cells ---\n$t = $d.Tables.Item(1)\n\n$cellRows = @(\n    @(\"14\u00f75=\", \"48\u00f75=\", \"59\u00f72=\", \"65\u00f76=\", \"79\u00f75=\"),\n    @(\"47\u00f77=\", \"36\u00f79=\", \"16\u00f72=\", \"37\u00f77=\", \"74\u00f72=\"),\n    @(\"40\u00f74=\", \"31\u00f76=\", \"41\u00f75=\", \"23\u00f74=\", \"93\u00f77=\"),\n    @(\"96\u00f79=\", \"51\u00f73=\", \"46\u00f74=\", \"22\u00f74=\", \"23\u00f74=\"),\n    @(\"56\u00f79=\", \"17\u00f78=\", \"16\u00f74=\", \"39\u00f75=\", \"90\u00f78=\")\n)\n\n# Non-blank data rows are every 4th row starting at 1 (1,5,9,13,17 \u2014 1-based).\nfor ($i = 0; $i -lt $cellRows.Count; $i++) {\n    $rowIndex = 1 + ($i * 4)\n    $values = $cellRows[$i]\n    for ($c = 1; $c -le $values.Count; $c++) {\n        $t.Cell($rowIndex, $c).Range.Text = $values[$c - 1]\n    }\n}\n"}
